$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 12
$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Cells.Item($row, 1).Value = 42619.893819444442
$ws.Cells.Item($row, 2).Value = 20
$ws.Cells.Item($row, 3).Value = 60
$ws.Cells.Item($row, 4).Value = 39
$ws.Cells.Item($row, 5).Value = 60
$ws.Cells.Item($row, 6).Value = 35
$ws.Cells.Item($row, 7).Value = 27870
$ws.Cells.Item($row, 8).Value = 13695
$ws.Cells.Item($row, 9).Value = 721
$ws.Cells.Item($row, 10).Value = 156
$ws.Cells.Item($row, 11).Value = 100
$ws.Cells.Item($row, 12).Value = 9
$ws.Cells.Item($row, 13).Value = 5
$ws.Cells.Item($row, 14).Value = "Named"
